# The presentation currently carries two DrawingML themes:
#   ppt/theme/theme1.xml  -> "Integral" ("Red Violet" colour scheme), used by
#                             the slide master / all slide layouts / slides.
#   ppt/theme/theme2.xml  -> "Office Theme", used only by the notes master.
#
# The authored edit swaps the two themes' content (theme1 becomes the
# "Office Theme" colours, theme2 becomes the former "Integral" colours).
# The font scheme and format scheme (fills/lines/effects) are byte-for-byte
# identical between the two themes, so the only substantive difference is
# the 12-colour <a:clrScheme>. We reproduce that by pushing the "Office
# Theme" palette onto the live theme through the standard PowerPoint
# ThemeColorScheme object (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -
# in that fixed index order).

function HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation

# Target "Office Theme" colour scheme (what theme2.xml currently holds),
# applied in the fixed ThemeColorScheme index order.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToVbaRgb $officeThemeColors[$i - 1]
}

$p.Save()
